# Add a new second slide using the "Title and Content" layout
# (PpSlideLayout ppLayoutText = 2, which maps to slideLayout2.xml:
#  a <p:ph type="title"/> placeholder + a generic <p:ph idx="1"/> placeholder).
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(2, 2)

# Set the title text of the new slide.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Second slide"
